# "finished Python and SQL" — refresh the sales-summary sheet with the
# final Python/SQL-computed figures and drop the now-obsolete header
# row + the blank spacer row above the tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old dashboard-title row (row 1) and the blank spacer row
# (row 2) that used to sit above the data tables; everything below
# shifts up by two rows.
$ws.Rows("1:2").Delete()

# The shared "Tổng doanh số" column header becomes simply "Doanh số"
# everywhere it is used (three mini-tables share the same label).
$ws.Range("C1").Value = "Doanh số"
$ws.Range("I1").Value = "Doanh số"
$ws.Range("K2").Value = "Doanh số"

# Refresh the "số đã bán" / "doanh số" figures for the brand table
# (columns B/C) with the finalized numbers.
$ws.Range("C2").Value  = 180426979000
$ws.Range("C3").Value  = 445230000
$ws.Range("C4").Value  = 1037043000
$ws.Range("C5").Value  = 797035000
$ws.Range("C6").Value  = 384880000
$ws.Range("C7").Value  = 207986240
$ws.Range("C8").Value  = 1953184000
$ws.Range("C9").Value  = 290433000
$ws.Range("C10").Value = 10192371140
$ws.Range("C11").Value = 233820000
$ws.Range("C12").Value = 4359920556
$ws.Range("C13").Value = 1291847000
$ws.Range("C14").Value = 160354260133
$ws.Range("C15").Value = 2494325860
$ws.Range("C16").Value = 179640000
$ws.Range("C18").Value = 1143591000
$ws.Range("C19").Value = 255598000
$ws.Range("C20").Value = 122184429081

# Refresh the "best-selling Chinese phone" mini-table (E/F columns).
$ws.Range("F2").Value = 24484

# Refresh the "revenue by country" mini-table (H/I columns).
$ws.Range("I2").Value = 181464022000
$ws.Range("I3").Value = 162307444133
$ws.Range("I4").Value = 140975520017

# Refresh the S23 vs. rest-of-Samsung ratio table (K/L/M columns).
$ws.Range("L2").Value = 47.61
$ws.Range("M2").Value = 52.39
$ws.Range("L3").Value = 82.23
$ws.Range("M3").Value = 17.77
